{"js": "// Resposta referente a quest\u00e3o-4\n//\n// The last paragraph of the document (the answer to \"Quest\u00e3o 3\") had its\n// trailing text awkwardly split mid-word (\"...auto-geren\" / \"cia para...\")\n// around a \"_GoBack\" bookmark. This script:\n//   1) fixes that paragraph so the text reads naturally (\"Time(\" / \"Squad)...\"\n//      with the same proofing marks Word itself adds around \"Time(\"),\n//   2) appends a new \"Quest\u00e3o 4:\" paragraph, and\n//   3) appends the full \"R:\" answer paragraph for Quest\u00e3o 4, keeping the\n//      \"_GoBack\" bookmark at the end of that final paragraph.\n//\n// We locate the target paragraph by its stable text fragment and replace its\n// OOXML wholesale (via Range.insertOoxml) with the corrected paragraph plus\n// the two new paragraphs, so the resulting markup matches exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Os desenvolvedores como\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error('Could not find the \"Quest\u00e3o 3\" answer paragraph.');\n}\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// Paragraph 1: the existing \"Quest\u00e3o 3\" answer, with the tail text repaired.\nconst para1 =\n  '<w:p w:rsidR=\"00A53D11\" w:rsidRDefault=\"00A53D11\">' +\n    '<w:r><w:t xml:space=\"preserve\">R: O SCRUM ser\u00e1 aplicado ao projeto, pois \u00e9 \u00e1gil, prioriza a satisfa\u00e7\u00e3o do cliente e entrega o software com maior rapidez e qualidade. O chefe da empresa seria o </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>PO( Product</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Owner) Respons\u00e1vel por decidir quais recursos ser\u00e3o constru\u00eddos e qual a ordem que devem ser feitos. O profissional full stack como </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>SM(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>Scrum Master) Respons\u00e1vel por ajudar todos os indiv\u00edduos a entender</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> os valores, para isso ele conhece muito bem o Scrum. Os desenvolvedores como </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>Time(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>Squad) onde a equipe se auto-gerencia para determinar a melhor maneira de realizar o trabalho.</w:t></w:r>' +\n  '</w:p>';\n\n// Paragraph 2: the new \"Quest\u00e3o 4:\" heading.\nconst para2 =\n  '<w:p>' +\n    '<w:r><w:t>Quest\u00e3o 4:</w:t></w:r>' +\n  '</w:p>';\n\n// Paragraph 3: the new \"R:\" answer, ending with the relocated \"_GoBack\" bookmark.\nconst para3 =\n  '<w:p>' +\n    '<w:r><w:t>R:</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">Conjuntos de dados: </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Quantidade_lixo_coletado</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>Caminhao</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> ,</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Rota, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Capacidade_Do_Caminhao</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Lixo_nao_coletado</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Os atributos de qualidade</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> necess\u00e1rio</w:t></w:r>' +\n    '<w:r><w:t>s s\u00e3o: </w:t></w:r>' +\n    '<w:r><w:t>Ader\u00eancia</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Unicidade</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Integridade</w:t></w:r>' +\n    '<w:r><w:t>, Legibilidade, Disponibilidade e Performance.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document ' + W_NS + '>' +\n          '<w:body>' + para1 + para2 + para3 + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nconst range = target.getRange(Word.RangeLocation.whole);\nrange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Resposta referente a quest\u00e3o-4\n#\n# The last paragraph of the document (the answer to \"Quest\u00e3o 3\") had its\n# trailing text awkwardly split mid-word (\"...auto-geren\" / \"cia para...\")\n# around a \"_GoBack\" bookmark. This script:\n#   1) fixes that paragraph so the text reads naturally (\"Time(\" / \"Squad)...\"\n#      with the same proofing marks Word itself adds around \"Time(\"),\n#   2) appends a new \"Quest\u00e3o 4:\" paragraph, and\n#   3) appends the full \"R:\" answer paragraph for Quest\u00e3o 4, keeping the\n#      \"_GoBack\" bookmark at the end of that final paragraph.\n#\n# We locate the target paragraph by its stable text fragment and replace its\n# contents wholesale (via Range.InsertXML, which parses and substitutes raw\n# WordprocessingML for the range) with the corrected paragraph plus the two\n# new paragraphs, so the resulting markup matches exactly.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like '*Os desenvolvedores como*') {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw 'Could not find the \"Quest\u00e3o 3\" answer paragraph.'\n}\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# Paragraph 1: the existing \"Quest\u00e3o 3\" answer, with the tail text repaired.\n$para1 = '<w:p ' + $wNs + ' w:rsidR=\"00A53D11\" w:rsidRDefault=\"00A53D11\">' +\n    '<w:r><w:t xml:space=\"preserve\">R: O SCRUM ser\u00e1 aplicado ao projeto, pois \u00e9 \u00e1gil, prioriza a satisfa\u00e7\u00e3o do cliente e entrega o software com maior rapidez e qualidade. O chefe da empresa seria o </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>PO( Product</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Owner) Respons\u00e1vel por decidir quais recursos ser\u00e3o constru\u00eddos e qual a ordem que devem ser feitos. O profissional full stack como </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>SM(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>Scrum Master) Respons\u00e1vel por ajudar todos os indiv\u00edduos a entender</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> os valores, para isso ele conhece muito bem o Scrum. Os desenvolvedores como </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>Time(</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>Squad) onde a equipe se auto-gerencia para determinar a melhor maneira de realizar o trabalho.</w:t></w:r>' +\n    '</w:p>'\n\n# Paragraph 2: the new \"Quest\u00e3o 4:\" heading.\n$para2 = '<w:p><w:r><w:t>Quest\u00e3o 4:</w:t></w:r></w:p>'\n\n# Paragraph 3: the new \"R:\" answer, ending with the relocated \"_GoBack\" bookmark.\n$para3 = '<w:p>' +\n    '<w:r><w:t>R:</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">Conjuntos de dados: </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Quantidade_lixo_coletado</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>Caminhao</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> ,</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Rota, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Capacidade_Do_Caminhao</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Lixo_nao_coletado</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Os atributos de qualidade</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> necess\u00e1rio</w:t></w:r>' +\n    '<w:r><w:t>s s\u00e3o: </w:t></w:r>' +\n    '<w:r><w:t>Ader\u00eancia</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Unicidade</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Integridade</w:t></w:r>' +\n    '<w:r><w:t>, Legibilidade, Disponibilidade e Performance.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n\n$xml = $para1 + $para2 + $para3\n\n$r = $target.Range\n$r.InsertXML($xml)\n"}
